$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "Test Case ID" values that keep their row position
#    (Log In, Search Page, Product Listing, Emailing Reports)
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,2).Value = "“1434”"
$ws.Cells.Item(3,2).Value = "“161 162 45”"
$ws.Cells.Item(4,2).Value = "“431 53 551 561 562”"
$ws.Cells.Item(7,2).Value = "“007”"

# ---------------------------------------------------------------------------
# 2. "COD Order" and "Checking Filters" swap places: COD Order used to be
#    row 6 and now becomes row 5; Checking Filters used to be row 5 and now
#    becomes row 6. Write the new contents for both rows directly (also
#    picks up the "Checking Filters" execution flag change NO -> YES and the
#    new / updated Test Case IDs) rather than trying to physically move
#    ranges around.
# ---------------------------------------------------------------------------

# New row 5: COD Order
$ws.Cells.Item(5,1).Value = "COD Order"
$ws.Cells.Item(5,2).Value = "“1434 431 53 551 561 562 `n612 8121 8123 8225 8471 8472 911”"
$ws.Cells.Item(5,3).Value = "YES"
$ws.Cells.Item(5,4).Value = "login"
$ws.Cells.Item(5,5).Value = "productCatalogPage"
$ws.Cells.Item(5,6).Value = "productDetailPage"
$ws.Cells.Item(5,7).Value = "cartCheck"
$ws.Cells.Item(5,8).Value = "checkout"

# New row 6: Checking Filters
$ws.Cells.Item(6,1).Value = "Checking Filters"
$ws.Cells.Item(6,2).Value = "“431 441 442 443 444 445 446`n447 448 45”"
$ws.Cells.Item(6,3).Value = "YES"
$ws.Cells.Item(6,4).Value = "productCatalogPage"
$ws.Cells.Item(6,5).Value = "applyFilters"
$ws.Cells.Item(6,6).Value = ""
$ws.Cells.Item(6,7).Value = ""
$ws.Cells.Item(6,8).Value = ""

# ---------------------------------------------------------------------------
# 3. The Test Case ID cells for these two rows now wrap their (longer,
#    multi-line) text, and the rows grow taller to fit it.
# ---------------------------------------------------------------------------
$ws.Cells.Item(5,2).WrapText = $true
$ws.Cells.Item(6,2).WrapText = $true
$ws.Rows.Item(5).RowHeight = 23.95
$ws.Rows.Item(6).RowHeight = 23.95

# ---------------------------------------------------------------------------
# 4. Column widths were tweaked slightly (column B grew a lot to host the
#    wrapped Test Case ID text, the rest shrank a touch).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.469387755102
$ws.Columns.Item(2).ColumnWidth = 26.9540816326531
$ws.Columns.Item(3).ColumnWidth = 12.5561224489796
$ws.Columns.Item(4).ColumnWidth = 21.734693877551
$ws.Columns.Item(5).ColumnWidth = 23.3520408163265
$ws.Columns.Item(6).ColumnWidth = 22.9489795918367
$ws.Columns.Item(7).ColumnWidth = 15.9285714285714
$ws.Columns.Item(8).ColumnWidth = 11.8775510204082

# ---------------------------------------------------------------------------
# 5. Active cell / selection moved from C7 to C4.
# ---------------------------------------------------------------------------
$ws.Range("C4").Select()
